$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (existing G->H, H->I), to make room
# for a new "6" model column between the existing "5" and "7" columns.
$ws.Columns.Item(7).Insert()

# Populate the newly inserted column G with the header and data values.
$ws.Range("G1").Value = 6
$ws.Range("G2").Value = 97.75466458084237
$ws.Range("G3").Value = 97.99027506983477
$ws.Range("G4").Value = 97.74054129477329
$ws.Range("G5").Value = 97.82963451003522
$ws.Range("G6").Value = 97.93273523665029
$ws.Range("G7").Value = 98.04907842579468
$ws.Range("G8").Value = 98.03280297904736
$ws.Range("G9").Value = 97.93975008160015
$ws.Range("G10").Value = 97.8580304454531
